# Remove api response formatting
# Appends a new data row (row 93) to each of the 4 worksheets, mirroring
# the structure of the existing rows (time, length, ID, actual length,
# checksum + their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = 45770.93995556713
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x5a"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 346
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = 45770.79742635417
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x5a"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 346
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = 45770.94318100694
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x5a"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 346
        I = 3
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = 45771.00693100694
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x5a"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 346
        I = 3
    }
)

foreach ($entry in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $newRow = 93

    $ws.Cells.Item($newRow, 1).Value = $entry.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $entry.B
    $ws.Cells.Item($newRow, 3).Value = $entry.C
    $ws.Cells.Item($newRow, 4).Value = $entry.D
    $ws.Cells.Item($newRow, 5).Value = $entry.E
    $ws.Cells.Item($newRow, 6).Value = $entry.F
    $ws.Cells.Item($newRow, 7).Value = $entry.G
    $ws.Cells.Item($newRow, 8).Value = $entry.H
    $ws.Cells.Item($newRow, 9).Value = $entry.I
}
